$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 119.0815153333333
$ws.Cells.Item(2, 8).Value = 357.244546
$ws.Cells.Item(2, 9).Value = 0.431812569872284
$ws.Cells.Item(2, 10).Value = 0.4318125698722839
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.563107
$ws.Cells.Item(2, 14).Value = 7.689321
$ws.Cells.Item(2, 15).Value = 0.09258576031692413
$ws.Cells.Item(2, 16).Value = 0.09258576031692412
$ws.Cells.Item(2, 17).Value = 305.218665521474
$ws.Cells.Item(2, 18).Value = 2746.967989693266
$ws.Cells.Item(2, 19).Value = 0.03997969509603034
$ws.Cells.Item(2, 20).Value = 0.03997969509603033

$ws.Cells.Item(3, 7).Value = 119.0815153333333
$ws.Cells.Item(3, 8).Value = 357.244546
$ws.Cells.Item(3, 9).Value = 0.431812569872284
$ws.Cells.Item(3, 10).Value = 0.4318125698722839
$ws.Cells.Item(3, 15).Value = 0.160803024221502
$ws.Cells.Item(3, 16).Value = 0.160803024221502
$ws.Cells.Item(3, 17).Value = 530.1040278407971
$ws.Cells.Item(3, 18).Value = 4770.936250567174
$ws.Cells.Item(3, 19).Value = 0.06943676713232193
$ws.Cells.Item(3, 20).Value = 0.06943676713232191

$ws.Cells.Item(4, 7).Value = 119.0815153333333
$ws.Cells.Item(4, 8).Value = 357.244546
$ws.Cells.Item(4, 9).Value = 0.431812569872284
$ws.Cells.Item(4, 10).Value = 0.4318125698722839
$ws.Cells.Item(4, 13).Value = 0.7887020000000001
$ws.Cells.Item(4, 14).Value = 2.366106
$ws.Cells.Item(4, 15).Value = 0.02848986575023154
$ws.Cells.Item(4, 16).Value = 0.02848986575023154
$ws.Cells.Item(4, 17).Value = 93.91982930643069
$ws.Cells.Item(4, 18).Value = 845.2784637578761
$ws.Cells.Item(4, 19).Value = 0.01230228214492385
$ws.Cells.Item(4, 20).Value = 0.01230228214492384

$ws.Cells.Item(5, 7).Value = 119.0815153333333
$ws.Cells.Item(5, 8).Value = 357.244546
$ws.Cells.Item(5, 9).Value = 0.431812569872284
$ws.Cells.Item(5, 10).Value = 0.4318125698722839
$ws.Cells.Item(5, 13).Value = 19.880183
$ws.Cells.Item(5, 14).Value = 59.640549
$ws.Cells.Item(5, 15).Value = 0.7181213497113423
$ws.Cells.Item(5, 16).Value = 0.7181213497113423
$ws.Cells.Item(5, 17).Value = 2367.362316743973
$ws.Cells.Item(5, 18).Value = 21306.26085069575
$ws.Cells.Item(5, 19).Value = 0.3100938254990079
$ws.Cells.Item(5, 20).Value = 0.3100938254990079

$ws.Cells.Item(6, 9).Value = 0.4460879372303943
$ws.Cells.Item(6, 10).Value = 0.4460879372303942
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.563107
$ws.Cells.Item(6, 14).Value = 7.689321
$ws.Cells.Item(6, 15).Value = 0.09258576031692413
$ws.Cells.Item(6, 16).Value = 0.09258576031692412
$ws.Cells.Item(6, 17).Value = 315.308942828964
$ws.Cells.Item(6, 18).Value = 2837.780485460676
$ws.Cells.Item(6, 19).Value = 0.04130139083668438
$ws.Cells.Item(6, 20).Value = 0.04130139083668437

$ws.Cells.Item(7, 9).Value = 0.4460879372303943
$ws.Cells.Item(7, 10).Value = 0.4460879372303942
$ws.Cells.Item(7, 15).Value = 0.160803024221502
$ws.Cells.Item(7, 16).Value = 0.160803024221502
$ws.Cells.Item(7, 19).Value = 0.07173228937537897
$ws.Cells.Item(7, 20).Value = 0.07173228937537897

$ws.Cells.Item(8, 9).Value = 0.4460879372303943
$ws.Cells.Item(8, 10).Value = 0.4460879372303942
$ws.Cells.Item(8, 13).Value = 0.7887020000000001
$ws.Cells.Item(8, 14).Value = 2.366106
$ws.Cells.Item(8, 15).Value = 0.02848986575023154
$ws.Cells.Item(8, 16).Value = 0.02848986575023154
$ws.Cells.Item(8, 17).Value = 97.02474138890402
$ws.Cells.Item(8, 18).Value = 873.2226725001361
$ws.Cells.Item(8, 19).Value = 0.01270898544449165
$ws.Cells.Item(8, 20).Value = 0.01270898544449164

$ws.Cells.Item(9, 9).Value = 0.4460879372303943
$ws.Cells.Item(9, 10).Value = 0.4460879372303942
$ws.Cells.Item(9, 13).Value = 19.880183
$ws.Cells.Item(9, 14).Value = 59.640549
$ws.Cells.Item(9, 15).Value = 0.7181213497113423
$ws.Cells.Item(9, 16).Value = 0.7181213497113423
$ws.Cells.Item(9, 17).Value = 2445.625362100116
$ws.Cells.Item(9, 18).Value = 22010.62825890104
$ws.Cells.Item(9, 19).Value = 0.3203452715738393
$ws.Cells.Item(9, 20).Value = 0.3203452715738393

$ws.Cells.Item(10, 7).Value = 33.50679633333333
$ws.Cells.Item(10, 8).Value = 100.520389
$ws.Cells.Item(10, 9).Value = 0.1215021138451521
$ws.Cells.Item(10, 10).Value = 0.121502113845152
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.563107
$ws.Cells.Item(10, 14).Value = 7.689321
$ws.Cells.Item(10, 15).Value = 0.09258576031692413
$ws.Cells.Item(10, 16).Value = 0.09258576031692412
$ws.Cells.Item(10, 17).Value = 85.88150422954101
$ws.Cells.Item(10, 18).Value = 772.9335380658689
$ws.Cells.Item(10, 19).Value = 0.01124936559046688
$ws.Cells.Item(10, 20).Value = 0.01124936559046687

$ws.Cells.Item(11, 7).Value = 33.50679633333333
$ws.Cells.Item(11, 8).Value = 100.520389
$ws.Cells.Item(11, 9).Value = 0.1215021138451521
$ws.Cells.Item(11, 10).Value = 0.121502113845152
$ws.Cells.Item(11, 15).Value = 0.160803024221502
$ws.Cells.Item(11, 16).Value = 0.160803024221502
$ws.Cells.Item(11, 17).Value = 149.1590667671768
$ws.Cells.Item(11, 18).Value = 1342.431600904591
$ws.Cells.Item(11, 19).Value = 0.01953790735560568
$ws.Cells.Item(11, 20).Value = 0.01953790735560568

$ws.Cells.Item(12, 7).Value = 33.50679633333333
$ws.Cells.Item(12, 8).Value = 100.520389
$ws.Cells.Item(12, 9).Value = 0.1215021138451521
$ws.Cells.Item(12, 10).Value = 0.121502113845152
$ws.Cells.Item(12, 13).Value = 0.7887020000000001
$ws.Cells.Item(12, 14).Value = 2.366106
$ws.Cells.Item(12, 15).Value = 0.02848986575023154
$ws.Cells.Item(12, 16).Value = 0.02848986575023154
$ws.Cells.Item(12, 17).Value = 26.42687728169267
$ws.Cells.Item(12, 18).Value = 237.841895535234
$ws.Cells.Item(12, 19).Value = 0.003461578911817731
$ws.Cells.Item(12, 20).Value = 0.00346157891181773

$ws.Cells.Item(13, 7).Value = 33.50679633333333
$ws.Cells.Item(13, 8).Value = 100.520389
$ws.Cells.Item(13, 9).Value = 0.1215021138451521
$ws.Cells.Item(13, 10).Value = 0.121502113845152
$ws.Cells.Item(13, 13).Value = 19.880183
$ws.Cells.Item(13, 14).Value = 59.640549
$ws.Cells.Item(13, 15).Value = 0.7181213497113423
$ws.Cells.Item(13, 16).Value = 0.7181213497113423
$ws.Cells.Item(13, 17).Value = 666.1212428503957
$ws.Cells.Item(13, 18).Value = 5995.09118565356
$ws.Cells.Item(13, 19).Value = 0.08725326198726178
$ws.Cells.Item(13, 20).Value = 0.08725326198726176

$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.16474
$ws.Cells.Item(14, 8).Value = 0.49422
$ws.Cells.Item(14, 9).Value = 0.000597379052169715
$ws.Cells.Item(14, 10).Value = 0.000597379052169715
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 2.563107
$ws.Cells.Item(14, 14).Value = 7.689321
$ws.Cells.Item(14, 15).Value = 0.09258576031692413
$ws.Cells.Item(14, 16).Value = 0.09258576031692412
$ws.Cells.Item(14, 17).Value = 0.42224624718
$ws.Cells.Item(14, 18).Value = 3.80021622462
$ws.Cells.Item(14, 19).Value = 0.00005530879374253655
$ws.Cells.Item(14, 20).Value = 0.00005530879374253654

$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.16474
$ws.Cells.Item(15, 8).Value = 0.49422
$ws.Cells.Item(15, 9).Value = 0.000597379052169715
$ws.Cells.Item(15, 10).Value = 0.000597379052169715
$ws.Cells.Item(15, 15).Value = 0.160803024221502
$ws.Cells.Item(15, 16).Value = 0.160803024221502
$ws.Cells.Item(15, 17).Value = 0.7333576273533333
$ws.Cells.Item(15, 18).Value = 6.600218646179999
$ws.Cells.Item(15, 19).Value = 0.00009606035819546461
$ws.Cells.Item(15, 20).Value = 0.00009606035819546461

$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.16474
$ws.Cells.Item(16, 8).Value = 0.49422
$ws.Cells.Item(16, 9).Value = 0.000597379052169715
$ws.Cells.Item(16, 10).Value = 0.000597379052169715
$ws.Cells.Item(16, 13).Value = 0.7887020000000001
$ws.Cells.Item(16, 14).Value = 2.366106
$ws.Cells.Item(16, 15).Value = 0.02848986575023154
$ws.Cells.Item(16, 16).Value = 0.02848986575023154
$ws.Cells.Item(16, 17).Value = 0.12993076748
$ws.Cells.Item(16, 18).Value = 1.16937690732
$ws.Cells.Item(16, 19).Value = 0.00001701924899831574
$ws.Cells.Item(16, 20).Value = 0.00001701924899831574

$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 0.6666666666666666
$ws.Cells.Item(17, 7).Value = 0.16474
$ws.Cells.Item(17, 8).Value = 0.49422
$ws.Cells.Item(17, 9).Value = 0.000597379052169715
$ws.Cells.Item(17, 10).Value = 0.000597379052169715
$ws.Cells.Item(17, 13).Value = 19.880183
$ws.Cells.Item(17, 14).Value = 59.640549
$ws.Cells.Item(17, 15).Value = 0.7181213497113423
$ws.Cells.Item(17, 16).Value = 0.7181213497113423
$ws.Cells.Item(17, 17).Value = 3.27506134742
$ws.Cells.Item(17, 18).Value = 29.47555212678
$ws.Cells.Item(17, 19).Value = 0.0004289906512333982
$ws.Cells.Item(17, 20).Value = 0.0004289906512333982
